$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text representation (avoid numeric auto-conversion)
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '51.867.27'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -1.12%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.921.99'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '357.94'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +1.25%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '109.68'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  -3.40%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.569'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +1.36%  '
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.630'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -0.01%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '39.27'
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -3.04%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0879'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +1.45%  '
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.70%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '19.61'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -2.98%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.90'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +0.05%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '3.385.71'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.07%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '2.928.57'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.06%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.988'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -1.07%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '51.884.81'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -1.08%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '3.37'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +0.78%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '7.60'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -1.20%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.09'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -2.89%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.0₃0983'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.45%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '71.05'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  -0.28%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '269.95'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.83%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.84'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +0.61%  '
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +13.24%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '26.99'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -0.43%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.62'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +16.65%  '
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +0.02%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.107'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +12.80%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '10.59'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -0.88%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '38.14'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -0.31%  '
$ws.Range('B33').NumberFormat = "@"
$ws.Range('B33').Value = 'RenderToken'
$ws.Range('C33').NumberFormat = "@"
$ws.Range('C33').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '6.06'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -2.23%  '
$ws.Range('B34').NumberFormat = "@"
$ws.Range('B34').Value = 'Toncoin'
$ws.Range('C34').NumberFormat = "@"
$ws.Range('C34').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.18'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -3.50%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '52.49'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -1.42%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.0444'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -2.47%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '1.00'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +0.11%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.25'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -3.15%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '18.36'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -3.61%  '
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -4.03%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -0.76%  '
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +1.91%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '23.01'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -5.27%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '119.29'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -2.65%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.19'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +0.09%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '3.49'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -2.93%  '
$ws.Range('B47').NumberFormat = "@"
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').NumberFormat = "@"
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '2.47'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -5.39%  '
$ws.Range('B48').NumberFormat = "@"
$ws.Range('B48').Value = 'Maker'
$ws.Range('C48').NumberFormat = "@"
$ws.Range('C48').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.128.17'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -4.22%  '
$ws.Range('B49').NumberFormat = "@"
$ws.Range('B49').Value = 'TheGraph'
$ws.Range('C49').NumberFormat = "@"
$ws.Range('C49').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.251'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -6.04%  '
$ws.Range('B50').NumberFormat = "@"
$ws.Range('B50').Value = 'BEAM'
$ws.Range('C50').NumberFormat = "@"
$ws.Range('C50').Value = 'https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0334'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -1.87%  '
$ws.Range('B51').NumberFormat = "@"
$ws.Range('B51').Value = 'FraxShare'
$ws.Range('C51').NumberFormat = "@"
$ws.Range('C51').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '9.16'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -0.10%  '
